$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.651186608805858
$ws.Range("C2").Value = 0.3382732131055093
$ws.Range("D2").Value = 0.5597078713389152
$ws.Range("E2").Value = 0.1943433289082286
$ws.Range("G2").Value = 2.062430422614398
$ws.Range("H2").Value = 1.689165945989373
$ws.Range("J2").Value = 0.07884750812304198
$ws.Range("M2").Value = 0.6770906075155807
$ws.Range("N2").Value = 1.911571362190486

$ws.Range("B3").Value = 1.549061444081417
$ws.Range("C3").Value = 0.3114674180820316
$ws.Range("D3").Value = 0.5564725585529686
$ws.Range("E3").Value = 0.1940779725982331
$ws.Range("G3").Value = 2.033396609915854
$ws.Range("H3").Value = 1.68286875242913
$ws.Range("J3").Value = 0.07916035238530128
$ws.Range("M3").Value = 0.6513079661504193
$ws.Range("N3").Value = 1.93147008805834

$ws.Range("B4").Value = 1.487287611493286
$ws.Range("C4").Value = 0.2951818148058578
$ws.Range("D4").Value = 0.5547536221390317
$ws.Range("E4").Value = 0.1939988157787305
$ws.Range("G4").Value = 2.016913980249967
$ws.Range("H4").Value = 1.679906984083317
$ws.Range("J4").Value = 0.07937824518145398
$ws.Range("M4").Value = 0.6358753592373745
$ws.Range("N4").Value = 1.944371355053633

$ws.Range("B5").Value = 1.462347365127641
$ws.Range("C5").Value = 0.2885885065411742
$ws.Range("D5").Value = 0.5541203494441618
$ws.Range("E5").Value = 0.193987609728179
$ws.Range("G5").Value = 2.010533473417524
$ws.Range("H5").Value = 1.678926909167046
$ws.Range("J5").Value = 0.07947352758191251
$ws.Range("M5").Value = 0.6296862020989806
$ws.Range("N5").Value = 1.949800168251087

$ws.Range("B6").Value = 1.458220102212636
$ws.Range("C6").Value = 0.2874962919598261
$ws.Range("D6").Value = 0.5540192521527132
$ws.Range("E6").Value = 0.1939870199982288
$ws.Range("G6").Value = 2.009494256681222
$ws.Range("H6").Value = 1.678777852530004
$ws.Range("J6").Value = 0.07948974110662377
$ws.Range("M6").Value = 0.6286645126785686
$ws.Range("N6").Value = 1.950711961719463

$ws.Range("B7").Value = 1.486950316302341
$ws.Range("C7").Value = 0.2950927206684923
$ws.Range("D7").Value = 0.5547448095713179
$ws.Range("E7").Value = 0.1939985794337886
$ws.Range("G7").Value = 2.016826571129428
$ws.Range("H7").Value = 1.679892848718453
$ws.Range("J7").Value = 0.07937950391786863
$ws.Range("M7").Value = 0.6357914866313408
$ws.Range("N7").Value = 1.94444387630822

$ws.Range("B8").Value = 1.615779814143423
$ws.Range("C8").Value = 0.3289943860352764
$ws.Range("D8").Value = 0.558536757117821
$ws.Range("E8").Value = 0.194234437520695
$ws.Range("G8").Value = 2.052139508773251
$ws.Range("H8").Value = 1.686806506300428
$ws.Range("J8").Value = 0.0789500211012264
$ws.Range("M8").Value = 0.6681179262112593
$ws.Range("N8").Value = 1.918290371995724

$ws.Range("B9").Value = 1.875871075697944
$ws.Range("C9").Value = 0.3968704406700567
$ws.Range("D9").Value = 0.5681001486779991
$ws.Range("E9").Value = 0.1953625840324378
$ws.Range("G9").Value = 2.132138566711262
$ws.Range("H9").Value = 1.707574508289412
$ws.Range("J9").Value = 0.07831257733798225
$ws.Range("M9").Value = 0.7346877889733179
$ws.Range("N9").Value = 1.872442815619962

$ws.Range("B10").Value = 2.071621067362116
$ws.Range("C10").Value = 0.4476230934918135
$ws.Range("D10").Value = 0.5764311560214139
$ws.Range("E10").Value = 0.1965989080633292
$ws.Range("G10").Value = 2.197594380403331
$ws.Range("H10").Value = 1.727276783131032
$ws.Range("J10").Value = 0.07796916026171274
$ws.Range("M10").Value = 0.7855674583161658
$ws.Range("N10").Value = 1.842097337966038

$ws.Range("B11").Value = 2.161711526190686
$ws.Range("C11").Value = 0.4709114998544237
$ws.Range("D11").Value = 0.5805063334875058
$ws.Range("E11").Value = 0.1972502302963832
$ws.Range("G11").Value = 2.228850787581479
$ws.Range("H11").Value = 1.737215751852489
$ws.Range("J11").Value = 0.07784007978664675
$ws.Range("M11").Value = 0.8091496472673469
$ws.Range("N11").Value = 1.829022348869508

$ws.Range("B12").Value = 2.19597813796895
$ws.Range("C12").Value = 0.4797596083274129
$ws.Range("D12").Value = 0.5820906590661821
$ws.Range("E12").Value = 0.1975096828493719
$ws.Range("G12").Value = 2.240901782126997
$ws.Range("H12").Value = 1.741120603507284
$ws.Range("J12").Value = 0.07779510534627221
$ws.Range("M12").Value = 0.8181429348386331
$ws.Range("N12").Value = 1.824176528295368

$ws.Range("B13").Value = 2.188591458459939
$ws.Range("C13").Value = 0.4778527001197403
$ws.Range("D13").Value = 0.5817476145091121
$ws.Range("E13").Value = 0.1974532349112792
$ws.Range("G13").Value = 2.238296795955307
$ws.Range("H13").Value = 1.740273332107876
$ws.Range("J13").Value = 0.07780461765298696
$ws.Range("M13").Value = 0.8162032507600259
$ws.Range("N13").Value = 1.825215466673285

$ws.Range("B14").Value = 2.164527620583669
$ws.Range("C14").Value = 0.4716388494954344
$ws.Range("D14").Value = 0.580635851669598
$ws.Range("E14").Value = 0.1972713187132662
$ws.Range("G14").Value = 2.22983791216285
$ws.Range("H14").Value = 1.737534172327088
$ws.Range("J14").Value = 0.07783630142353459
$ws.Range("M14").Value = 0.8098882609060638
$ws.Range("N14").Value = 1.82862156504082

$ws.Range("B15").Value = 2.149807570991811
$ws.Range("C15").Value = 0.4678365143207088
$ws.Range("D15").Value = 0.5799602263626014
$ws.Range("E15").Value = 0.1971615589300058
$ws.Range("G15").Value = 2.224684644522029
$ws.Range("H15").Value = 1.735874767302079
$ws.Range("J15").Value = 0.07785621734588943
$ws.Range("M15").Value = 0.8060283935432011
$ws.Range("N15").Value = 1.830721639591836

$ws.Range("B16").Value = 2.06575454300031
$ws.Range("C16").Value = 0.4461052173797952
$ws.Range("D16").Value = 0.5761705850669046
$ws.Range("E16").Value = 0.1965581339183515
$ws.Range("G16").Value = 2.195581658413687
$ws.Range("H16").Value = 1.726646960973568
$ws.Range("J16").Value = 0.0779781423502115
$ws.Range("M16").Value = 0.7840351311713931
$ws.Range("N16").Value = 1.842966539181162

$ws.Range("B17").Value = 2.014458773608112
$ws.Range("C17").Value = 0.4328254910337819
$ws.Range("D17").Value = 0.5739189273820102
$ws.Range("E17").Value = 0.1962107425558948
$ws.Range("G17").Value = 2.178108467223524
$ws.Range("H17").Value = 1.721236612265329
$ws.Range("J17").Value = 0.07805989278869419
$ws.Range("M17").Value = 0.7706551380115485
$ws.Range("N17").Value = 1.850665551208955

$ws.Range("B18").Value = 1.985052874265648
$ws.Range("C18").Value = 0.4252062151333007
$ws.Range("D18").Value = 0.5726506860980294
$ws.Range("E18").Value = 0.1960192994813603
$ws.Range("G18").Value = 2.168197553777304
$ws.Range("H18").Value = 1.718216569254707
$ws.Range("J18").Value = 0.07810946801127017
$ws.Range("M18").Value = 0.7630004052516739
$ws.Range("N18").Value = 1.855162434906667

$ws.Range("B19").Value = 1.975113344990518
$ws.Range("C19").Value = 0.4226296876708489
$ws.Range("D19").Value = 0.5722258900796646
$ws.Range("E19").Value = 0.1959559164700053
$ws.Range("G19").Value = 2.164865736191985
$ws.Range("H19").Value = 1.717209785495953
$ws.Range("J19").Value = 0.07812669199791245
$ws.Range("M19").Value = 0.7604156882039632
$ws.Range("N19").Value = 1.856696771584289

$ws.Range("B20").Value = 2.019909136200624
$ws.Range("C20").Value = 0.4342371844109607
$ws.Range("D20").Value = 0.5741558402562816
$ws.Range("E20").Value = 0.1962468568000588
$ws.Range("G20").Value = 2.17995409665744
$ws.Range("H20").Value = 1.721803041698877
$ws.Range("J20").Value = 0.07805092591270224
$ws.Range("M20").Value = 0.7720752059562841
$ws.Range("N20").Value = 1.849838873054772

$ws.Range("B21").Value = 2.171591635576419
$ws.Range("C21").Value = 0.4734632095788811
$ws.Range("D21").Value = 0.5809612858245998
$ws.Range("E21").Value = 0.1973244040159088
$ws.Range("G21").Value = 2.232316643619725
$ws.Range("H21").Value = 1.738334892014279
$ws.Range("J21").Value = 0.07782688911517255
$ws.Range("M21").Value = 0.8117414072291069
$ws.Range("N21").Value = 1.827618246843443

$ws.Range("B22").Value = 2.271607525083994
$ws.Range("C22").Value = 0.4992706388929378
$ws.Range("D22").Value = 0.5856489124038831
$ws.Range("E22").Value = 0.1981033281190392
$ws.Range("G22").Value = 2.267791986272414
$ws.Range("H22").Value = 1.749962692576105
$ws.Range("J22").Value = 0.07770323394790069
$ws.Range("M22").Value = 0.8380342850147571
$ws.Range("N22").Value = 1.813710347095089

$ws.Range("B23").Value = 2.21814594617058
$ws.Range("C23").Value = 0.4854809431089393
$ws.Range("D23").Value = 0.5831250519739797
$ws.Range("E23").Value = 0.1976807592942507
$ws.Range("G23").Value = 2.248742769355175
$ws.Range("H23").Value = 1.743681122058007
$ws.Range("J23").Value = 0.07776714703722476
$ws.Range("M23").Value = 0.8239674023279946
$ws.Range("N23").Value = 1.821076848627449

$ws.Range("B24").Value = 2.017444763485116
$ws.Range("C24").Value = 0.4335989098506161
$ws.Range("D24").Value = 0.5740486501391615
$ws.Range("E24").Value = 0.1962305037716163
$ws.Range("G24").Value = 2.179119268316128
$ws.Range("H24").Value = 1.721546677315757
$ws.Range("J24").Value = 0.07805497181434262
$ws.Range("M24").Value = 0.7714330761431256
$ws.Range("N24").Value = 1.850212394324771

$ws.Range("B25").Value = 1.80469941938594
$ws.Range("C25").Value = 0.3783553781983926
$ws.Range("D25").Value = 0.5652844467021794
$ws.Range("E25").Value = 0.1949859595523549
$ws.Range("G25").Value = 2.109332783032102
$ws.Range("H25").Value = 1.701179342361058
$ws.Range("J25").Value = 0.07846309337674029
$ws.Range("M25").Value = 0.716335246645464
$ws.Range("N25").Value = 1.884261152909367

Write-Host "Updated pl_mw values for 380 kV case"
